$wb = $excel.ActiveWorkbook

$wsSurvey = $wb.Worksheets.Item("survey")
$wsPromptTypes = $wb.Worksheets.Item("prompt_types")
$wsModel = $wb.Worksheets.Item("model")

# --- prompt_types: add new "adate" prompt type row (row 4) ---
$wsPromptTypes.Cells.Item(4, 1).Value = "adate"
$wsPromptTypes.Cells.Item(4, 2).Value = "string"
$wsPromptTypes.Cells.Item(4, 3).Value = "string"
$wsPromptTypes.Cells.Item(4, 4).Value = "Save only mm.dd.yyyy with support for ?? at all positions"

# --- survey: switch elementType from custom_date to adate ---
$wsSurvey.Cells.Item(19, 4).Value = "adate"
$wsSurvey.Cells.Item(21, 4).Value = "adate"
$wsSurvey.Cells.Item(47, 4).Value = "adate"

# --- model: switch elementType from custom_date to adate ---
$wsModel.Cells.Item(3, 2).Value = "adate"
$wsModel.Cells.Item(10, 2).Value = "adate"
$wsModel.Cells.Item(11, 2).Value = "adate"

# model row 9 (GRVISITDATA) becomes async_assign_date, matching the
# formatting (wrap text style + taller row) already used by the
# prompt_types!A3 "async_assign_date" cell.
$wsPromptTypes.Range("A3").Copy() | Out-Null
$wsModel.Cells.Item(9, 2).PasteSpecial(-4163) | Out-Null

# --- selections, to mirror where the editor left the cursor ---
[void]$wsSurvey.Range("D4").Select()
[void]$wsPromptTypes.Range("D9").Select()
[void]$wsModel.Range("F11").Select()
